$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column E ("Abonas").
#    This shifts the existing E (Carga Horaria) -> F, F (Banco do Dia) -> G,
#    G (Saldo Acumulado) -> H. Styles/formatting of the shifted cells move
#    along automatically with the insert.
$ws.Columns("E:E").Insert()

# 2. Set the new header for column E.
$ws.Range("E1").Value = "Abonas"

# 3. Fill the new "Abonas" column with "00:00" for every existing data row
#    (rows 2 through 66).
$ws.Range("E2:E66").Value = "00:00"

# 4. Update row 66 ("Marcacoes" and "Horas Trabalhadas" for 13/01/2026).
$ws.Range("C66").Value = "07:21 | 12:00 | 13:00 | 16:21"
$ws.Range("D66").Value = "08:00"

# 5. Append a brand-new row 67 for 14/01/2026.
#    Temporarily force column A to Text format so the ISO date string is not
#    auto-converted into a date serial number, then clear the formatting
#    again so the cell ends up with the default (no explicit) style, same
#    as the value, but still holding the literal text.
$ws.Range("A67").NumberFormat = "@"
$ws.Range("A67").Value = "2026-01-14"
$ws.Range("A67").ClearFormats()
$ws.Range("B67").Value = "14/01/2026"
$ws.Range("C67").Value = "08:00 | 11:34 | 12:34 | 16:05"
$ws.Range("D67").Value = "07:05"
$ws.Range("E67").Value = "00:00"
$ws.Range("F67").Value = "08:00"
$ws.Range("G67").Value = "-00:55"
$ws.Range("H67").Value = "07:33"
